$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Borders.Item(8).LineStyle = 1
$ws.Range("C2").Borders.Item(8).Weight = 2
Write-Output "done"
